{"js": "const pairs = [\n  [\"65\u00d747=3055\", \"11\u00d743=473\"],\n  [\"45\u00d789=4005\", \"43\u00d725=1075\"],\n  [\"36\u00d796=3456\", \"75\u00d791=6825\"],\n  [\"20\u00d728=560\", \"43\u00d716=688\"],\n  [\"45\u00d763=2835\", \"72\u00d712=864\"],\n  [\"18\u00d773=1314\", \"50\u00d722=1100\"],\n  [\"59\u00d781=4779\", \"75\u00d779=5925\"],\n  [\"61\u00d755=3355\", \"87\u00d742=3654\"],\n  [\"71\u00d745=3195\", \"97\u00d755=5335\"],\n  [\"71\u00d779=5609\", \"11\u00d750=550\"],\n  [\"13\u00d718=234\", \"53\u00d775=3975\"],\n  [\"37\u00d711=407\", \"51\u00d772=3672\"],\n  [\"17\u00d720=340\", \"37\u00d721=777\"],\n  [\"28\u00d793=2604\", \"47\u00d773=3431\"],\n  [\"36\u00d713=468\", \"78\u00d736=2808\"],\n  [\"84\u00d752=4368\", \"81\u00d773=5913\"],\n  [\"94\u00d712=1128\", \"48\u00d740=1920\"],\n  [\"86\u00d737=3182\", \"75\u00d760=4500\"],\n  [\"97\u00d750=4850\", \"77\u00d779=6083\"],\n  [\"41\u00d770=2870\", \"74\u00d730=2220\"],\n  [\"79\u00d771=5609\", \"91\u00d765=5915\"],\n  [\"97\u00d769=6693\", \"68\u00d747=3196\"],\n  [\"41\u00d762=2542\", \"69\u00d795=6555\"],\n  [\"89\u00d766=5874\", \"94\u00d785=7990\"],\n  [\"46\u00d727=1242\", \"96\u00d765=6240\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @('65\u00d747=3055', '11\u00d743=473'),\n    @('45\u00d789=4005', '43\u00d725=1075'),\n    @('36\u00d796=3456', '75\u00d791=6825'),\n    @('20\u00d728=560', '43\u00d716=688'),\n    @('45\u00d763=2835', '72\u00d712=864'),\n    @('18\u00d773=1314', '50\u00d722=1100'),\n    @('59\u00d781=4779', '75\u00d779=5925'),\n    @('61\u00d755=3355', '87\u00d742=3654'),\n    @('71\u00d745=3195', '97\u00d755=5335'),\n    @('71\u00d779=5609', '11\u00d750=550'),\n    @('13\u00d718=234', '53\u00d775=3975'),\n    @('37\u00d711=407', '51\u00d772=3672'),\n    @('17\u00d720=340', '37\u00d721=777'),\n    @('28\u00d793=2604', '47\u00d773=3431'),\n    @('36\u00d713=468', '78\u00d736=2808'),\n    @('84\u00d752=4368', '81\u00d773=5913'),\n    @('94\u00d712=1128', '48\u00d740=1920'),\n    @('86\u00d737=3182', '75\u00d760=4500'),\n    @('97\u00d750=4850', '77\u00d779=6083'),\n    @('41\u00d770=2870', '74\u00d730=2220'),\n    @('79\u00d771=5609', '91\u00d765=5915'),\n    @('97\u00d769=6693', '68\u00d747=3196'),\n    @('41\u00d762=2542', '69\u00d795=6555'),\n    @('89\u00d766=5874', '94\u00d785=7990'),\n    @('46\u00d727=1242', '96\u00d765=6240')\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
